$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 480-481, pushing the existing rows 480:505
# down to 482:507 (dimension grows from A1:T505 to A1:T507).
$ws.Rows("480:481").Insert()

# New row 480: same dimension/product metadata as the row that used to sit
# there, with updated date / variety / volume / price / unit / origin data.
$ws.Range("A480").Value = 9
$ws.Range("B480").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C480").Value = "Metropolitana"
$ws.Range("D480").Value = "2021-11-16"
$ws.Range("E480").Value = 13
$ws.Range("F480").Value = "Fruta"
$ws.Range("G480").Value = 100102
$ws.Range("H480").Value = "Cítricos"
$ws.Range("I480").Value = 100102005
$ws.Range("J480").Value = "Naranja"
$ws.Range("K480").Value = "New Hall"
$ws.Range("L480").Value = "Primera"
$ws.Range("M480").Value = 500
$ws.Range("N480").Value = 7500
$ws.Range("O480").Value = 8000
$ws.Range("P480").Value = 7780
$ws.Range("Q480").Value = "$/malla 18 kilos"
$ws.Range("R480").Value = "Región de O'Higgins"
$ws.Range("S480").Value = 432
$ws.Range("T480").Value = 18

# New row 481
$ws.Range("A481").Value = 9
$ws.Range("B481").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C481").Value = "Metropolitana"
$ws.Range("D481").Value = "2021-11-16"
$ws.Range("E481").Value = 13
$ws.Range("F481").Value = "Fruta"
$ws.Range("G481").Value = 100102
$ws.Range("H481").Value = "Cítricos"
$ws.Range("I481").Value = 100102005
$ws.Range("J481").Value = "Naranja"
$ws.Range("K481").Value = "Valencia"
$ws.Range("L481").Value = "Primera"
$ws.Range("M481").Value = 450
$ws.Range("N481").Value = 8000
$ws.Range("O481").Value = 8000
$ws.Range("P481").Value = 8000
$ws.Range("Q481").Value = "$/caja 15 kilos granel"
$ws.Range("R481").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S481").Value = 533
$ws.Range("T481").Value = 15
